# Modified and tested the code for merging database.
# Swap the contents of rows 3 and 7 (A/B columns) on the active sheet,
# then leave the selection on A3:B3 (active cell A3), matching the
# edit made in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values of the two rows being swapped.
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$a7 = $ws.Range("A7").Value2
$b7 = $ws.Range("B7").Value2

# Swap row 3 and row 7 contents.
$ws.Range("A3").Value = $a7
$ws.Range("B3").Value = $b7
$ws.Range("A7").Value = $a3
$ws.Range("B7").Value = $b3

# Update the visible selection to A3:B3 (active cell A3).
$ws.Range("A3:B3").Select()
